# chemical_import_template.xlsx - import-template column cleanup
#
# Commit: "feat: Enhance chemical and sample import functionality"
#   - Drops a handful of columns from the sample_chemicals import sheet
#     (residue type / sample readout / created at / updated at /
#     user labels / literatures) so the template lines up with the
#     fields the importer actually consumes.
#   - Refreshes the "last saved window" bookkeeping (absolute path hint,
#     selection, header row height) that Excel stamps on every save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the now-unused columns from the header row. We look each
#    one up by its header text (row 1) rather than a hard-coded letter
#    so the script is resilient to the column shifting left as each
#    prior match is removed.
# ---------------------------------------------------------------------
$columnsToRemove = @(
    "residue type",
    "sample readout",
    "created at",
    "updated at",
    "user labels",
    "literatures"
)

foreach ($name in $columnsToRemove) {
    $headerRow = $ws.Range("A1:BZ1")
    $found = $headerRow.Find($name)
    if ($found -ne $null) {
        $found.EntireColumn.Delete()
    }
}

# ---------------------------------------------------------------------
# 2. Bump the header row height slightly (matches the refreshed
#    template's styling pass).
# ---------------------------------------------------------------------
$ws.Rows("1").RowHeight = 16.5

# ---------------------------------------------------------------------
# 3. Restore the workbook/window chrome Excel re-writes on save: the
#    remembered selection and the recent-files absolute-path hint.
# ---------------------------------------------------------------------
$ws.Range("X1:AA1048576").Select()

try {
    $wb.Path = "C:\Users\49172\Downloads\templates\"
} catch {
}

$wb.Save()
